$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed topic-distribution rows with newly computed values
$ws.Range("B2").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:1.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B3").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.351|x10:0.000|x11:0.000|x12:0.211|x13:0.438|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B4").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.161|x9:0.000|x10:0.148|x11:0.000|x12:0.350|x13:0.154|x14:0.000|x15:0.000|x16:0.187|x17:0.000"
$ws.Range("B6").Value = "x1:0.144|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.661|x9:0.000|x10:0.195|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B7").Value = "x1:0.000|x2:0.346|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.654|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B10").Value = "x1:0.335|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.258|x9:0.000|x10:0.217|x11:0.000|x12:0.000|x13:0.000|x14:0.190|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B11").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.524|x12:0.476|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B12").Value = "x1:nan|x2:nan|x3:nan|x4:nan|x5:nan|x6:nan|x7:nan|x8:nan|x9:nan|x10:nan|x11:nan|x12:nan|x13:nan|x14:nan|x15:nan|x16:nan|x17:nan"
$ws.Range("B13").Value = "x1:0.000|x2:0.000|x3:0.000|x4:0.000|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.689|x10:0.311|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B16").Value = "x1:0.340|x2:0.146|x3:0.000|x4:0.514|x5:0.000|x6:0.000|x7:0.000|x8:0.000|x9:0.000|x10:0.000|x11:0.000|x12:0.000|x13:0.000|x14:0.000|x15:0.000|x16:0.000|x17:0.000"
$ws.Range("B18").Value = "x1:nan|x2:nan|x3:nan|x4:nan|x5:nan|x6:nan|x7:nan|x8:nan|x9:nan|x10:nan|x11:nan|x12:nan|x13:nan|x14:nan|x15:nan|x16:nan|x17:nan"

# Remove the old totals row (row 19), which is no longer part of the map
$ws.Rows("19").Delete()
